$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old "Is Significant" column (F) needs to shift right to make room for
# two new columns: Observed (F) and Expected (G).
$ws.Range("F:G").Insert()

# Header cells for the two new columns (F, G)
$ws.Range("F1").Value = "Observed"
$ws.Range("G1").Value = "Expected"

# Observed / Expected count-string pairs per feature row (rows 2..24)
$observed = @(
    "[505  84] ; [ 76 142]",
    "[485 117] ; [107 118]",
    "[492 105] ; [133  92]",
    "[505  88] ; [141  81]",
    "[531  71] ; [159  66]",
    "[566  31] ; [182  41]",
    "[561  35] ; [182  41]",
    "[299 290] ; [ 65 156]",
    "[383 208] ; [103 119]",
    "[240 358] ; [ 51 171]",
    "[425 167] ; [123 100]",
    "[345 234] ; [ 93 129]",
    "[329 264] ; [ 85 139]",
    "[559  32] ; [188  33]",
    "[331 246] ; [ 89 130]",
    "[446 153] ; [134  89]",
    "[304 281] ; [ 81 141]",
    "[284 311] ; [ 72 150]",
    "[188 389] ; [ 41 176]",
    "[262 314] ; [ 65 146]",
    "[467 130] ; [148  75]",
    "[422 178] ; [129  95]",
    "[481 120] ; [154  70]"
)

$expected = @(
    "[424.05080545 164.94919455] ; [156.94919455  61.05080545]",
    "[430.93591294 171.06408706] ; [161.06408706  63.93591294]",
    "[453.92335766 143.07664234] ; [171.07664234  53.92335766]",
    "[470.03435583 122.96564417] ; [175.96564417  46.03435583]",
    "[502.2732769  99.7267231] ; [187.7267231  37.2732769]",
    "[544.5804878  52.4195122] ; [203.4195122  19.5804878]",
    "[540.69352869  55.30647131] ; [202.30647131  20.69352869]",
    "[264.68641975 324.31358025] ; [ 99.31358025 121.68641975]",
    "[353.29151292 237.70848708] ; [132.70848708  89.29151292]",
    "[212.21707317 385.78292683] ; [ 78.78292683 143.21707317]",
    "[398.05644172 193.94355828] ; [149.94355828  73.05644172]",
    "[316.60674157 262.39325843] ; [121.39325843 100.60674157]",
    "[300.49204406 292.50795594] ; [113.50795594 110.49204406]",
    "[543.6908867  47.3091133] ; [203.3091133  17.6908867]",
    "[304.44723618 272.55276382] ; [115.55276382 103.44723618]",
    "[422.65206813 176.34793187] ; [157.34793187  65.65206813]",
    "[279.08921933 305.91078067] ; [105.91078067 116.08921933]",
    "[259.26560588 335.73439412] ; [ 96.73439412 125.26560588]",
    "[166.41435768 410.58564232] ; [ 62.58564232 154.41435768]",
    "[239.32909784 336.67090216] ; [ 87.67090216 123.32909784]",
    "[447.75 149.25] ; [167.25  55.75]",
    "[401.21359223 198.78640777] ; [149.78640777  74.21359223]",
    "[462.58787879 138.41212121] ; [172.41212121  51.58787879]"
)

for ($i = 0; $i -lt 23; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $observed[$i]
}
for ($i = 0; $i -lt 23; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $expected[$i]
}

